$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cells in row 3 ---
$ws.Range("Q3").Value = 42630.89999985695
$ws.Range("R3").Value = 274999482515.4157
$ws.Range("T3").Value = 5927549283187.865
$ws.Range("AA3").Value = 'userData\20250802233910-1'

# --- Update existing cells in row 4 ---
$ws.Range("A4").Value = 208936824749.1004
$ws.Range("B4").Value = 560558279.3021399
$ws.Range("C4").Value = 2838405.422452337
$ws.Range("D4").Value = 205285332057.718
$ws.Range("E4").Value = 3088096006.65783
$ws.Range("F4").Value = 205845890337.0201
$ws.Range("G4").Value = 3090934412.080283
$ws.Range("M4").Value = 208936824749.1004
$ws.Range("Q4").Value = 0.382000207901001
$ws.Range("R4").Value = 208884473638.7114
$ws.Range("S4").Value = 208936824749.1004
$ws.Range("T4").Value = 52351110.38903809
$ws.Range("Z4").Value = -1
$ws.Range("AA4").Value = 'userData\20250826141509-1'

# --- Add new rows 5-9 ---
# Row 5
$ws.Range("A5").Value = 208936824749.1004
$ws.Range("B5").Value = 560558279.3021399
$ws.Range("C5").Value = 2838405.422452337
$ws.Range("D5").Value = 205285332057.718
$ws.Range("E5").Value = 3088096006.65783
$ws.Range("F5").Value = 205845890337.0201
$ws.Range("G5").Value = 3090934412.080283
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 208936824749.1004
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 0.3420000076293945
$ws.Range("R5").Value = 208884473638.7114
$ws.Range("S5").Value = 208936824749.1004
$ws.Range("T5").Value = 52351110.38903809
$ws.Range("U5").Value = 'optimal'
$ws.Range("V5").Value = 'costs_emissionlimit'
$ws.Range("W5").Value = -1
$ws.Range("X5").Value = -1
$ws.Range("Y5").Value = 1
$ws.Range("Z5").Value = -1
$ws.Range("AA5").Value = 'userData\20250826142033-1'
$ws.Range("AB5").Value = 0
$ws.Range("AC5").Value = 0
$ws.Range("AD5").Value = 0
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 0

# Row 6
$ws.Range("A6").Value = 212761070453.3561
$ws.Range("B6").Value = 431929795.0159013
$ws.Range("C6").Value = 483280.7941501038
$ws.Range("D6").Value = 210117737617.5461
$ws.Range("E6").Value = 2210919760
$ws.Range("F6").Value = 210549667412.562
$ws.Range("G6").Value = 2211403040.79415
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 212761070453.3561
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 0.2049999237060547
$ws.Range("R6").Value = 208884473638.7115
$ws.Range("S6").Value = 212761070453.3561
$ws.Range("T6").Value = 3876596814.644623
$ws.Range("U6").Value = 'optimal'
$ws.Range("V6").Value = 'costs_emissionlimit'
$ws.Range("W6").Value = -1
$ws.Range("X6").Value = -1
$ws.Range("Y6").Value = 1
$ws.Range("Z6").Value = -1
$ws.Range("AA6").Value = 'userData\20250826142349-1'
$ws.Range("AB6").Value = 0
$ws.Range("AC6").Value = 0
$ws.Range("AD6").Value = 0
$ws.Range("AE6").Value = 0
$ws.Range("AF6").Value = 0

# Row 7
$ws.Range("A7").Value = 212761070453.3561
$ws.Range("B7").Value = 431929795.0159013
$ws.Range("C7").Value = 483280.7941501038
$ws.Range("D7").Value = 210117737617.5461
$ws.Range("E7").Value = 2210919760
$ws.Range("F7").Value = 210549667412.562
$ws.Range("G7").Value = 2211403040.79415
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 212761070453.3561
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 0.2089998722076416
$ws.Range("R7").Value = 208884473638.7115
$ws.Range("S7").Value = 212761070453.3561
$ws.Range("T7").Value = 3876596814.644623
$ws.Range("U7").Value = 'optimal'
$ws.Range("V7").Value = 'costs_emissionlimit'
$ws.Range("W7").Value = -1
$ws.Range("X7").Value = -1
$ws.Range("Y7").Value = 1
$ws.Range("Z7").Value = -1
$ws.Range("AA7").Value = 'userData\20250826142803-1'
$ws.Range("AB7").Value = 0
$ws.Range("AC7").Value = 0
$ws.Range("AD7").Value = 0
$ws.Range("AE7").Value = 0
$ws.Range("AF7").Value = 0

# Row 8
$ws.Range("A8").Value = 212555257606.3859
$ws.Range("B8").Value = 391248723.5784202
$ws.Range("C8").Value = 563827.5931751186
$ws.Range("D8").Value = 209922525295.2143
$ws.Range("E8").Value = 2240919760
$ws.Range("F8").Value = 210313774018.7927
$ws.Range("G8").Value = 2241483587.593174
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 212555257606.3859
$ws.Range("N8").Value = 0
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = 0
$ws.Range("Q8").Value = 0.1640000343322754
$ws.Range("R8").Value = 208897880768.9483
$ws.Range("S8").Value = 212555257606.3859
$ws.Range("T8").Value = 3657376837.437622
$ws.Range("U8").Value = 'optimal'
$ws.Range("V8").Value = 'costs_emissionlimit'
$ws.Range("W8").Value = -1
$ws.Range("X8").Value = -1
$ws.Range("Y8").Value = 1
$ws.Range("Z8").Value = -1
$ws.Range("AA8").Value = 'userData\20250826144826-1'
$ws.Range("AB8").Value = 0
$ws.Range("AC8").Value = 0
$ws.Range("AD8").Value = 0
$ws.Range("AE8").Value = 0
$ws.Range("AF8").Value = 0

# Row 9
$ws.Range("A9").Value = 327302522268.1901
$ws.Range("B9").Value = 100493544796.2502
$ws.Range("C9").Value = 6257650125.644318
$ws.Range("D9").Value = 212401019790.5789
$ws.Range("E9").Value = 8150307555.716727
$ws.Range("F9").Value = 312894564586.8291
$ws.Range("G9").Value = 14407957681.36105
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 327302522268.1901
$ws.Range("N9").Value = 0
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = 0
$ws.Range("Q9").Value = 23106.82899999619
$ws.Range("R9").Value = 324700548326.465
$ws.Range("S9").Value = 327302522268.1901
$ws.Range("T9").Value = 2601973941.725159
$ws.Range("U9").Value = 'optimal'
$ws.Range("V9").Value = 'costs_emissionlimit'
$ws.Range("W9").Value = -1
$ws.Range("X9").Value = -1
$ws.Range("Y9").Value = 1
$ws.Range("Z9").Value = -1
$ws.Range("AA9").Value = 'userData\20250829000334-1'
$ws.Range("AB9").Value = 0
$ws.Range("AC9").Value = 0
$ws.Range("AD9").Value = 0
$ws.Range("AE9").Value = 0
$ws.Range("AF9").Value = 0

